$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 (board numbers above the existing grid)
$ws.Range("F4").Value = 14
$ws.Range("H4").Value = 15
$ws.Range("I4").Value = 16
$ws.Range("K4").Value = 18
$ws.Range("L4").Value = 19
$ws.Range("N4").Value = 21

# New column D values (left side of the board)
$ws.Range("D6").Value = 13
$ws.Range("D7").Value = 12
$ws.Range("D9").Value = 11
$ws.Range("D11").Value = 9
$ws.Range("D12").Value = 8
$ws.Range("D14").Value = 6

# New column P values (right side of the board)
$ws.Range("P6").Value = 22
$ws.Range("P7").Value = 23
$ws.Range("P9").Value = 24
$ws.Range("P12").Value = 26
$ws.Range("P14").Value = 27

# New row 16 (board numbers below the existing grid)
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 4
$ws.Range("I16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("N16").Value = 0

# Update the view: the new selection moves to Q14
$ws.Range("Q14").Select()
